{"js": "// Remove the trailing \"Ver no Jupiter ...\" and \"\u00a9 2020 ...\" paragraphs\n// (and the blank paragraph that separated them from the requirement line),\n// leaving the requirement paragraph (\"LOM3049: ...\") followed directly by\n// the single blank paragraph that used to sit after the copyright line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the three paragraphs to drop by their text content so the script\n// does not depend on brittle absolute indices.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter' / copyright paragraphs to remove.\");\n}\n\n// The blank paragraph immediately preceding the \"Ver no Jupiter\" paragraph\n// (separating it from the requirement text) is removed too.\nlet blankIdx = jupiterIdx - 1;\nif (blankIdx < 0 || items[blankIdx].text !== \"\") {\n  blankIdx = -1;\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIdx].delete();\nitems[jupiterIdx].delete();\nif (blankIdx !== -1) {\n  items[blankIdx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" and \"\u00a9 2020 ...\" paragraphs\n# (and the blank paragraph that separated them from the requirement line),\n# leaving the requirement paragraph (\"LOM3049: ...\") followed directly by\n# the single blank paragraph that used to sit after the copyright line.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n# Locate the paragraphs to drop by their text content so the script does\n# not depend on brittle absolute indices.\n$jupiterIdx = -1\n$copyrightIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($jupiterIdx -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Powered by Jekyll*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -eq -1 -or $copyrightIdx -eq -1) {\n    throw \"Could not locate the 'Ver no Jupiter' / copyright paragraphs to remove.\"\n}\n\n# The blank paragraph immediately preceding the \"Ver no Jupiter\" paragraph\n# (separating it from the requirement text) is removed too.\n$blankIdx = $jupiterIdx - 1\n$blankText = $paras.Item($blankIdx).Range.Text.TrimEnd([char]13, [char]7)\nif ($blankIdx -lt 1 -or $blankText.Length -ne 0) {\n    $blankIdx = $jupiterIdx\n}\n\n$startPara = $paras.Item($blankIdx)\n$endPara = $paras.Item($copyrightIdx)\n\n$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$rng.Delete()\n"}
